$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Script1")

# Update the PathParam value in B3 to the new path parameter string
$ws.Range("B3").Value = "/api/users/{user}"

# Adjust column B width to fit the new content
$ws.Columns.Item(2).ColumnWidth = 16.42578125
